# Instruction.docx edit: replace the "Live cricket score webpage" project
# instructions with the "Mood Checker" project instructions, and merge the
# trailing bookmark paragraph into the "Primary programming language"
# paragraph, leaving one new empty paragraph at the end.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-ParagraphRuns {
    param($Index, $FirstText, $SecondText)

    $para = $d.Paragraphs($Index).Range
    # Include the paragraph mark so the whole paragraph's contents are
    # replaced (and its own <w:p> properties regenerated) in one shot.
    $rng = $d.Range($para.Start, $para.End)

    $xml = '<w:p ' + $wNs + '>' +
        '<w:r><w:t xml:space="preserve">' + $FirstText + '</w:t></w:r>' +
        '<w:r><w:t>' + $SecondText + '</w:t></w:r>' +
        '</w:p>'

    $null = $rng.InsertXML($xml)
}

# Paragraph 1: Project Name
Set-ParagraphRuns 1 "Project Name: " "Mood Checker"

# Paragraph 2: Project Description
Set-ParagraphRuns 2 "Project Description: " "Python AI application which takes input from user related to their mood and give them motivational or positive quote depending upon their mood"

# Paragraph 3: Problem statement
Set-ParagraphRuns 3 "Problem statement: App " "will help user to understand and boost their mood"

# Paragraph 4: Users
Set-ParagraphRuns 4 "Users: " "User will receive motivational quote as per their mood"

# Paragraph 5: Platforms
Set-ParagraphRuns 5 "Platforms: Platform should be " "web app"

# Paragraphs 6 & 7: merge the bookmark-only paragraph into the "Primary
# programming language" paragraph, then leave a fresh empty paragraph
# behind (where the bookmark paragraph used to be).
$p6 = $d.Paragraphs(6).Range
$p7 = $d.Paragraphs(7).Range
$mergeRng = $d.Range($p6.Start, $p7.End)

$mergeXml = '<w:p ' + $wNs + ' w:rsidR="00490496" w:rsidRDefault="00490496" w:rsidP="00490496">' +
    '<w:r><w:t>Primary programming language: Python Flask</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '</w:p>' +
    '<w:p ' + $wNs + '/>'

$null = $mergeRng.InsertXML($mergeXml)
